$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").Value = $null
$ws.Range("H10").Value = 7950
$ws.Range("J10").Value = 7950
$ws.Range("L10").Value = 7950
$ws.Range("N10").Value = -8536
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").Value = $null
$ws.Range("H93").Value = 29811.766
$ws.Range("J93").Value = 29811.766
$ws.Range("L93").Value = 29811.766
$ws.Range("N93").Value = -34803.766

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").Value = $null
$ws.Range("H9").Value = 60008
$ws.Range("I9").Value = 60008
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 60008
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -59838
$ws.Range("N9").Value = $null
$ws.Range("H14").Value = 3366.3333
$ws.Range("I14").Value = 149.5
$ws.Range("J14").Value = 9800
$ws.Range("K14").Value = 149.5
$ws.Range("L14").Value = 9800
$ws.Range("M14").Value = 25.5
$ws.Range("N14").Value = -10150
$ws.Range("H20").Value = 60008
$ws.Range("I20").Value = 60008
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 60008
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -59738
$ws.Range("N20").Value = $null
$ws.Range("H21").Value = 15000
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 15000
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 15000
$ws.Range("M21").Value = $null
$ws.Range("N21").Value = -15748
$ws.Range("H23").Value = 20000
$ws.Range("J23").Value = 20000
$ws.Range("L23").Value = 20000
$ws.Range("N23").Value = -20518
$ws.Range("H40").Value = 8725
$ws.Range("J40").Value = 8725
$ws.Range("L40").Value = 8725
$ws.Range("N40").Value = -9077
$ws.Range("H42").Value = 15000
$ws.Range("J42").Value = 15000
$ws.Range("L42").Value = 15000
$ws.Range("N42").Value = -15972

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 400
$ws.Range("I8").Value = 400
$ws.Range("K8").Value = 400
$ws.Range("M8").Value = -260
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").Value = $null
$ws.Range("H44").Value = 20000
$ws.Range("J44").Value = 20000
$ws.Range("L44").Value = 20000
$ws.Range("N44").Value = -20994

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 10125
$ws.Range("I38").Value = 2000
$ws.Range("J38").Value = 18250
$ws.Range("K38").Value = 2000
$ws.Range("L38").Value = 18250
$ws.Range("M38").Value = -1623
$ws.Range("N38").Value = -19004
$ws.Range("H46").Value = 10125
$ws.Range("I46").Value = 2000
$ws.Range("J46").Value = 18250
$ws.Range("K46").Value = 2000
$ws.Range("L46").Value = 18250
$ws.Range("M46").Value = -1789
$ws.Range("N46").Value = -18672

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5483.3335
$ws.Range("I3").Value = 3406.25
$ws.Range("J3").Value = 7857.143
$ws.Range("K3").Value = 10218.75
$ws.Range("L3").Value = 23571.429
$ws.Range("M3").Value = -10106.75
$ws.Range("N3").Value = -23795.429
$ws.Range("H106").Value = 4664.706
$ws.Range("J106").Value = 4664.706
$ws.Range("L106").Value = 13994.118
$ws.Range("N106").Value = -15886.118
$ws.Range("H122").Value = 697.2105
$ws.Range("I122").Value = 467.58066
$ws.Range("J122").Value = 1714.1428
$ws.Range("K122").Value = 4208.22594
$ws.Range("L122").Value = 15427.2852
$ws.Range("M122").Value = -1758.22594
$ws.Range("N122").Value = -20327.2852
$ws.Range("H131").Value = 920.47
$ws.Range("J131").Value = 951.56384
$ws.Range("L131").Value = 2854.69152
$ws.Range("N131").Value = -12934.69152

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 5500000
$ws.Range("I14").Value = 5500000
$ws.Range("K14").Value = 5500000
$ws.Range("M14").Value = -5499832
$ws.Range("H38").Value = 8012
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 8012
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 8012
$ws.Range("M38").Value = $null
$ws.Range("N38").Value = -8938

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1802
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").Value = $null
$ws.Range("H16").Value = 527.8946999999999
$ws.Range("I16").Value = 555.3333
$ws.Range("J16").Value = 425
$ws.Range("K16").Value = 555.3333
$ws.Range("L16").Value = 425
$ws.Range("M16").Value = -385.3333
$ws.Range("N16").Value = -765
$ws.Range("H35").Value = 3265.25
$ws.Range("I35").Value = 530.5
$ws.Range("J35").Value = 6000
$ws.Range("K35").Value = 530.5
$ws.Range("L35").Value = 6000
$ws.Range("M35").Value = -194.5
$ws.Range("N35").Value = -6672
$ws.Range("H47").Value = 20000
$ws.Range("J47").Value = 20000
$ws.Range("L47").Value = 20000
$ws.Range("N47").Value = -20980
$ws.Range("H52").Value = 20000
$ws.Range("J52").Value = 20000
$ws.Range("L52").Value = 20000
$ws.Range("N52").Value = -20466
$ws.Range("H126").Value = 1802
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").Value = $null
$ws.Range("H133").Value = 34125.8
$ws.Range("J133").Value = 34125.8
$ws.Range("L133").Value = 34125.8
$ws.Range("N133").Value = -39185.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 3702.5
$ws.Range("I10").Value = 405
$ws.Range("J10").Value = 7000
$ws.Range("K10").Value = 405
$ws.Range("L10").Value = 7000
$ws.Range("M10").Value = -236
$ws.Range("N10").Value = -7338
$ws.Range("H15").Value = 9000
$ws.Range("J15").Value = 9000
$ws.Range("L15").Value = 9000
$ws.Range("N15").Value = -9576
$ws.Range("H64").Value = 26500
$ws.Range("J64").Value = 26500
$ws.Range("L64").Value = 26500
$ws.Range("N64").Value = -26996
$ws.Range("H67").Value = 26500
$ws.Range("J67").Value = 26500
$ws.Range("L67").Value = 26500
$ws.Range("N67").Value = -28216
$ws.Range("H138").Value = 32143.2
$ws.Range("J138").Value = 32143.2
$ws.Range("L138").Value = 32143.2
$ws.Range("N138").Value = -42423.2
